# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) held a placeholder/old "Strike#" stat; this
# recomputes/rewrites it with the freshly simulated strikeout (K) values
# for each game row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value for column G
$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 1
    32 = 1
    33 = 2
    34 = 0
    35 = 0
    36 = 1
    37 = 2
    38 = 2
    39 = 2
    40 = 5
    41 = 2
    42 = 5
    43 = 0
    44 = 1
    45 = 3
    46 = 1
    48 = 1
    49 = 2
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 2
    55 = 1
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 0
    61 = 1
    62 = 0
    63 = 1
    66 = 3
    67 = 1
    68 = 2
    69 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
